$wb = $excel.ActiveWorkbook

# --- Rename existing "Sheet1" to "Veg" (keeps sheetId=2) ---
$veg = $wb.Worksheets.Item(2)
$veg.Name = "Veg"

# --- Insert a brand-new blank sheet "NonVeg" right after "Veg" (becomes sheetId=3) ---
$nonVeg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $veg)
$nonVeg.Name = "NonVeg"

# --- Populate "Veg" with the data-driven menu table:
#     header row (restaurants A2B / SRR / SRV) typed across first,
#     then each restaurant's Idly/Dosa/Vada/Poori/MasalDosa prices filled in column by column ---
$veg.Cells.Item(1, 1).Value = "A2B"
$veg.Cells.Item(1, 2).Value = "SRR"
$veg.Cells.Item(1, 3).Value = "SRV"

$colA = @("Idly-8",  "Dosa-30", "Vada-45", "Poori-60", "MasalDosa -75")
$colB = @("Idly-9",  "Dosa-40", "Vada-5",  "Poori-20", "MasalDosa -55")
$colC = @("Idly-5",  "Dosa-31", "Vada-40", "Poori-10", "MasalDosa -65")

for ($i = 0; $i -lt $colA.Length; $i++) { $veg.Cells.Item($i + 2, 1).Value = $colA[$i] }
for ($i = 0; $i -lt $colB.Length; $i++) { $veg.Cells.Item($i + 2, 2).Value = $colB[$i] }
for ($i = 0; $i -lt $colC.Length; $i++) { $veg.Cells.Item($i + 2, 3).Value = $colC[$i] }

# Header row (row 1) is bold
$veg.Range("A1:C1").Font.Bold = $true

# Auto-fit the three columns to their content
$veg.Columns("A:C").AutoFit()

# Portrait page orientation
$veg.PageSetup.Orientation = 1

# "Veg" becomes the active/selected sheet, with C13 selected
$null = $veg.Range("C13").Select()

Write-Host "Veg + NonVeg sheets created"
